# Applies the "cryptos list" price/volume refresh described by the commit
# "Updated cryptos list on Tue Oct  3 02:06:44 UTC 2023 with GitHub Actions".
# Rows 2-51 hold one cryptocurrency each (columns: A=rank, B=Coin, C=Link,
# D=Price, E=Volume(1h)). This update refreshes the Price/Volume figures
# that moved, and additionally re-ranks Stellar above BinanceUSD (rows 28/29
# swap their Coin/Link/Price values, each keeping its own refreshed volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.531.66'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '1.666.25'
$ws.Range("E3").Value = '  -3.39%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '215.14'
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("D6").Value = '0.514'
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '23.65'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").Value = '0.262'
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("D10").Value = '0.0622'
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("E11").Value = '  -2.61%  '
$ws.Range("D12").Value = '1.901.93'
$ws.Range("E12").Value = '  -3.38%  '
$ws.Range("D13").Value = '1.670.72'
$ws.Range("E13").Value = '  -3.30%  '
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").Value = '0.554'
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("D16").Value = '66.33'
$ws.Range("D17").Value = '248.18'
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("D18").Value = '27.561.96'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("D19").Value = '0.0₃0734'
$ws.Range("E19").Value = '  -2.97%  '
$ws.Range("D20").Value = '7.57'
$ws.Range("E20").Value = '  -4.37%  '
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = '4.50'
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("E23").Value = '  -4.31%  '
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  -5.72%  '
$ws.Range("D25").Value = '146.54'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").Value = '16.62'
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("E27").Value = '  -5.06%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '0.112'
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("B29").Value = 'BinanceUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").Value = '1.24'
$ws.Range("E30").Value = '  +4.01%  '
$ws.Range("D31").Value = '0.0511'
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").Value = '3.36'
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").Value = '1.470.29'
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("E34").Value = '  -5.40%  '
$ws.Range("E35").Value = '  -5.55%  '
$ws.Range("D36").Value = '0.942'
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("E37").Value = '  -1.09%  '
$ws.Range("D38").Value = '0.577'
$ws.Range("E38").Value = '  -6.30%  '
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("D40").Value = '69.90'
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -6.48%  '
$ws.Range("D43").Value = '5.45'
$ws.Range("E43").Value = '  -6.74%  '
$ws.Range("D44").Value = '2.22'
$ws.Range("E44").Value = '  -3.11%  '
$ws.Range("D45").Value = '1.810.36'
$ws.Range("E45").Value = '  -3.30%  '
$ws.Range("D46").Value = '0.792'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").Value = '1.70'
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("D48").Value = '89.67'
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("D50").Value = '42.19'
$ws.Range("E50").Value = '  +20.33%  '
$ws.Range("E51").Value = '  -3.12%  '

